$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.713998556137085
$ws.Range("B1").Value = 2.267333984375
$ws.Range("C1").Value = 2.440818309783936
$ws.Range("D1").Value = 3.147013425827026
$ws.Range("E1").Value = 1.808821558952332
